$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 109.083336
$ws.Range("I2").Value = 75.666664
$ws.Range("J2").Value = 209.33333
$ws.Range("K2").Value = 75.666664
$ws.Range("L2").Value = 209.33333
$ws.Range("M2").Value = 37.333336
$ws.Range("N2").Value = -435.33333
# row 11
$ws.Range("H11").Value = 749.7273
$ws.Range("I11").Value = 749.7273
$ws.Range("K11").Value = 749.7273
$ws.Range("M11").Value = -609.7273
# row 19
$ws.Range("H19").Value = 4637.125
$ws.Range("I19").Value = 4399
$ws.Range("K19").Value = 4399
$ws.Range("M19").Value = -4224
# row 30
$ws.Range("H30").Value = 750
$ws.Range("J30").Value = 750
$ws.Range("L30").Value = 2250
$ws.Range("N30").Value = -2452
# row 57
$ws.Range("H57").Value = 136464.6
$ws.Range("J57").Value = 136464.6
$ws.Range("L57").Value = 409393.8
$ws.Range("N57").Value = -410391.8
# row 113
$ws.Range("H113").Value = 1001
$ws.Range("I113").Value = 1001
$ws.Range("K113").Value = 1001
$ws.Range("M113").Value = 2253
# row 119
$ws.Range("H119").Value = 2000
$ws.Range("J119").Value = 2000
$ws.Range("L119").Value = 6000
$ws.Range("N119").Value = -15676
# row 138
$ws.Range("H138").Value = 1554.2162
$ws.Range("I138").Value = 1250.7812
$ws.Range("K138").Value = 3752.3436
$ws.Range("M138").Value = 1387.6564

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4910.109
$ws.Range("I32").Value = 3522.94
$ws.Range("K32").Value = 3522.94
$ws.Range("M32").Value = -3235.94
# row 61
$ws.Range("H61").Value = 13162
$ws.Range("I61").Value = 12550.25
$ws.Range("J61").Value = 14997.25
$ws.Range("K61").Value = 12550.25
$ws.Range("L61").Value = 14997.25
$ws.Range("M61").Value = -12338.25
$ws.Range("N61").Value = -15421.25
# row 74
$ws.Range("H74").Value = 6965.5557
$ws.Range("I74").Value = 7211.5415
$ws.Range("K74").Value = 7211.5415
$ws.Range("M74").Value = -6337.5415
# row 77
$ws.Range("H77").Value = 6965.5557
$ws.Range("I77").Value = 7211.5415
$ws.Range("K77").Value = 36057.7075
$ws.Range("M77").Value = -31689.7075
# row 122
$ws.Range("H122").Value = 5192.7646
$ws.Range("I122").Value = 4611.8
$ws.Range("J122").Value = 9550
$ws.Range("K122").Value = 13835.4
$ws.Range("L122").Value = 28650
$ws.Range("M122").Value = -11385.4
$ws.Range("N122").Value = -33550
# row 132
$ws.Range("H132").Value = 3187
$ws.Range("I132").Value = 3187
$ws.Range("K132").Value = 9561
$ws.Range("M132").Value = -7031
# row 136
$ws.Range("H136").Value = 13162
$ws.Range("I136").Value = 12550.25
$ws.Range("J136").Value = 14997.25
$ws.Range("K136").Value = 37650.75
$ws.Range("L136").Value = 44991.75
$ws.Range("M136").Value = -35100.75
$ws.Range("N136").Value = -50091.75
# row 139
$ws.Range("H139").Value = 69999.75
$ws.Range("J139").Value = 69999.75
$ws.Range("L139").Value = 69999.75
$ws.Range("N139").Value = -80279.75

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 266
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# row 70
$ws.Range("H70").Value = 169909
$ws.Range("J70").Value = 169909
$ws.Range("L70").Value = 169909
$ws.Range("N70").Value = -170495
# row 73
$ws.Range("H73").Value = 169909
$ws.Range("J73").Value = 169909
$ws.Range("L73").Value = 169909
$ws.Range("N73").Value = -171937
# row 80
$ws.Range("H80").Value = 2084388
$ws.Range("J80").Value = 2778790.8
$ws.Range("L80").Value = 2778790.8
$ws.Range("N80").Value = -2780786.8
# row 83
$ws.Range("H83").Value = 2084388
$ws.Range("J83").Value = 2778790.8
$ws.Range("L83").Value = 13893954
$ws.Range("N83").Value = -13903938
# row 132
$ws.Range("H132").Value = 101242.4
$ws.Range("J132").Value = 101242.4
$ws.Range("L132").Value = 101242.4
$ws.Range("N132").Value = -111362.4
# row 134
$ws.Range("H134").Value = 6325.533
$ws.Range("I134").Value = 6933.913
$ws.Range("K134").Value = 20801.739
$ws.Range("M134").Value = -18266.739

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2770.182
$ws.Range("I31").Value = 2546.973
$ws.Range("K31").Value = 2546.973
$ws.Range("M31").Value = -2251.973
# row 34
$ws.Range("H34").Value = 2770.182
$ws.Range("I34").Value = 2546.973
$ws.Range("K34").Value = 2546.973
$ws.Range("M34").Value = -2344.973
# row 122
$ws.Range("H122").Value = 3516.7693
$ws.Range("I122").Value = 2840.6667
$ws.Range("K122").Value = 8522.000100000001
$ws.Range("M122").Value = -6072.000100000001

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 915.6667
$ws.Range("I5").Value = 905.125
$ws.Range("K5").Value = 2715.375
$ws.Range("M5").Value = -2603.375
# row 37
$ws.Range("H37").Value = 95191.91
$ws.Range("J37").Value = 95191.91
$ws.Range("L37").Value = 285575.73
$ws.Range("N37").Value = -285799.73
# row 38
$ws.Range("H38").Value = 102.375
$ws.Range("I38").Value = 83.84614999999999
$ws.Range("J38").Value = 182.66667
$ws.Range("K38").Value = 251.53845
$ws.Range("L38").Value = 548.00001
$ws.Range("M38").Value = 95.46155000000002
$ws.Range("N38").Value = -1242.00001
# row 104
$ws.Range("H104").Value = 17564.1
$ws.Range("I104").Value = 11995.5
$ws.Range("J104").Value = 19950.643
$ws.Range("K104").Value = 35986.5
$ws.Range("L104").Value = 59851.929
$ws.Range("M104").Value = -33365.5
$ws.Range("N104").Value = -65093.929
# row 113
$ws.Range("H113").Value = 360.4
$ws.Range("I113").Value = 366.25
$ws.Range("J113").Value = 337
$ws.Range("K113").Value = 1098.75
$ws.Range("L113").Value = 1011
$ws.Range("M113").Value = 1071.25
$ws.Range("N113").Value = -5351
# row 128
$ws.Range("H128").Value = 184326.67
$ws.Range("I128").Value = 184326.67
$ws.Range("K128").Value = 552980.01
$ws.Range("M128").Value = -548000.01
# row 135
$ws.Range("H135").Value = 915.6667
$ws.Range("I135").Value = 905.125
$ws.Range("K135").Value = 8146.125
$ws.Range("M135").Value = -5611.125

$ws = $wb.Worksheets.Item("GSM")
# row 57
$ws.Range("H57").Value = 13809
$ws.Range("I57").Value = 3015
$ws.Range("K57").Value = 3015
$ws.Range("M57").Value = -2195
# row 102
$ws.Range("H102").Value = 4974.963
$ws.Range("I102").Value = 5805.5293
$ws.Range("K102").Value = 5805.5293
$ws.Range("M102").Value = -4183.5293
# row 136
$ws.Range("H136").Value = 53360.11
$ws.Range("J136").Value = 53360.11
$ws.Range("L136").Value = 160080.33
$ws.Range("N136").Value = -165180.33
# row 140
$ws.Range("H140").Value = 69796.38
$ws.Range("J140").Value = 69796.38
$ws.Range("L140").Value = 69796.38
$ws.Range("N140").Value = -80156.38

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1614.7693
$ws.Range("J22").Value = 2004
$ws.Range("L22").Value = 2004
$ws.Range("N22").Value = -2594
# row 27
$ws.Range("H27").Value = 1614.7693
$ws.Range("J27").Value = 2004
$ws.Range("L27").Value = 2004
$ws.Range("N27").Value = -2218
# row 40
$ws.Range("H40").Value = 11202.119
$ws.Range("I40").Value = 11038.23
$ws.Range("J40").Value = 13332.667
$ws.Range("K40").Value = 11038.23
$ws.Range("L40").Value = 13332.667
$ws.Range("M40").Value = -10902.23
$ws.Range("N40").Value = -13604.667
# row 99
$ws.Range("H99").Value = 59999
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
# row 136
$ws.Range("H136").Value = 3279501.5
$ws.Range("I136").Value = 3755106
$ws.Range("K136").Value = 11265318
$ws.Range("M136").Value = -11262768
# row 137
$ws.Range("H137").Value = 164747.75
$ws.Range("J137").Value = 179996
$ws.Range("L137").Value = 179996
$ws.Range("N137").Value = -190196

$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 331.66666
$ws.Range("I113").Value = 304.625
$ws.Range("J113").Value = 548
$ws.Range("K113").Value = 913.875
$ws.Range("L113").Value = 1644
$ws.Range("M113").Value = 1256.125
$ws.Range("N113").Value = -5984
# row 126
$ws.Range("H126").Value = 7018.0605
$ws.Range("I126").Value = 6846.85
$ws.Range("K126").Value = 20540.55
$ws.Range("M126").Value = -18070.55
